# Apply updates described in the commit:
# "Fixed the hit box not activating issue, changed charge states from bools to an enum
#  the excel has been updated to match."
#
# This adds a new "Status" column (D) with progress notes on a few tasks,
# and fills in a couple of previously-empty "Completerer of task" (B) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header
$ws.Range("D1").Value2 = "Status"

# Status notes for specific tasks
$ws.Range("D2").Value2 = "complete, the hitboxes are probably too big now that they actually work but that can be changed easily"
$ws.Range("D5").Value2 = "complete"
$ws.Range("D7").Value2 = "these work now I think, just need to get good feel"

# Fill in previously blank "Completerer of task" cells
$ws.Range("B12").Value2 = "           "
$ws.Range("B14").Value2 = "Harris"

# Autofit the new column to match its content width
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Update selection to reflect last-edited cell
$ws.Range("D17").Select() | Out-Null
